# 2022.7.14 21.37 for textures
# Adds two new block entries (red coral / sea shroom) to the "方块" (Blocks)
# sheet, mirrors the print setup tweak picked up for that sheet, and leaves
# the workbook focused on the sheet/cell that was being edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("方块")

# New row 11: 红珊瑚 (red coral)
$ws.Range("A11").Value = "红珊瑚"
$ws.Range("B11").Value = "四周需要有完整方块依附，否则会脱落"
$ws.Range("C11").Value = "有四个方向"
$ws.Range("D11").Value = "red_coral"

# New row 12: 海蘑菇 (sea shroom)
# (write D before B so new shared-string indices land in the same order
# Excel produced: 海蘑菇, sea_shroom, 下方方块为砂砾)
$ws.Range("A12").Value = "海蘑菇"
$ws.Range("D12").Value = "sea_shroom"
$ws.Range("B12").Value = "下方方块为砂砾"

# Pick up the print setup that now accompanies this sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Focus moved to this sheet/cell for the session that made the edit.
$ws.Activate() | Out-Null
$ws.Range("D12").Select() | Out-Null
